# Apply numeric corrections to the Malboro Profits workbook (scheduled runner sync).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1554
$ws.Range("I9").Value = 1500
$ws.Range("K9").Value = 1500
$ws.Range("M9").Value = -1331
$ws.Range("H19").Value = 50134.777
$ws.Range("I19").Value = 824
$ws.Range("K19").Value = 824
$ws.Range("M19").Value = -649
$ws.Range("H51").Value = 8625
$ws.Range("I51").Value = 8833.333000000001
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 8833.333000000001
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = -8349.333000000001
$ws.Range("N51").Value = -8968
$ws.Range("H106").Value = 19126
$ws.Range("I106").Value = 8835
$ws.Range("K106").Value = 8835
$ws.Range("M106").Value = -8204
$ws.Range("H113").Value = 18767.8
$ws.Range("J113").Value = 3643.5
$ws.Range("L113").Value = 3643.5
$ws.Range("N113").Value = -10151.5
$ws.Range("H137").Value = 14308.806
$ws.Range("I137").Value = 2688.5
$ws.Range("J137").Value = 23605.05
$ws.Range("K137").Value = 8065.5
$ws.Range("L137").Value = 70815.14999999999
$ws.Range("M137").Value = -5515.5
$ws.Range("N137").Value = -75915.14999999999
$ws.Range("H138").Value = 4174.5415
$ws.Range("J138").Value = 3982.9119
$ws.Range("L138").Value = 11948.7357
$ws.Range("N138").Value = -22228.7357

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4221.3716
$ws.Range("I32").Value = 2245.762
$ws.Range("K32").Value = 2245.762
$ws.Range("M32").Value = -1958.762
$ws.Range("H45").Value = 1799.55
$ws.Range("I45").Value = 1568.1333
$ws.Range("J45").Value = 2493.8
$ws.Range("K45").Value = 1568.1333
$ws.Range("L45").Value = 2493.8
$ws.Range("M45").Value = -1191.1333
$ws.Range("N45").Value = -3247.8
$ws.Range("H61").Value = 1184700.6
$ws.Range("I61").Value = 4108.9375
$ws.Range("J61").Value = 3283530.2
$ws.Range("K61").Value = 4108.9375
$ws.Range("L61").Value = 3283530.2
$ws.Range("M61").Value = -3896.9375
$ws.Range("N61").Value = -3283954.2
$ws.Range("H74").Value = 29507.4
$ws.Range("I74").Value = 11954.889
$ws.Range("K74").Value = 11954.889
$ws.Range("M74").Value = -11080.889
$ws.Range("H77").Value = 29507.4
$ws.Range("I77").Value = 11954.889
$ws.Range("K77").Value = 59774.44499999999
$ws.Range("M77").Value = -55406.44499999999
$ws.Range("H136").Value = 1184700.6
$ws.Range("I136").Value = 4108.9375
$ws.Range("J136").Value = 3283530.2
$ws.Range("K136").Value = 12326.8125
$ws.Range("L136").Value = 9850590.600000001
$ws.Range("M136").Value = -9776.8125
$ws.Range("N136").Value = -9855690.600000001
$ws.Range("H137").Value = 95833.336
$ws.Range("J137").Value = 95833.336
$ws.Range("L137").Value = 95833.336
$ws.Range("N137").Value = -106033.336

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1687.1666
$ws.Range("I105").Value = 1567.8182
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1567.8182
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 179.1818000000001
$ws.Range("N105").Value = -6494

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13814.091
$ws.Range("I16").Value = 5133.3335
$ws.Range("K16").Value = 5133.3335
$ws.Range("M16").Value = -4846.3335
$ws.Range("H113").Value = 13814.091
$ws.Range("I113").Value = 5133.3335
$ws.Range("K113").Value = 5133.3335
$ws.Range("M113").Value = -2963.3335

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 879.2222
$ws.Range("I44").Value = 428.4
$ws.Range("J44").Value = 1442.75
$ws.Range("K44").Value = 1285.2
$ws.Range("L44").Value = 4328.25
$ws.Range("M44").Value = -887.1999999999998
$ws.Range("N44").Value = -5124.25
$ws.Range("H68").Value = 3776.1667
$ws.Range("I68").Value = 1200
$ws.Range("J68").Value = 4010.3635
$ws.Range("K68").Value = 3600
$ws.Range("L68").Value = 12031.0905
$ws.Range("M68").Value = -2789
$ws.Range("N68").Value = -13653.0905
$ws.Range("H71").Value = 3776.1667
$ws.Range("I71").Value = 1200
$ws.Range("J71").Value = 4010.3635
$ws.Range("K71").Value = 10800
$ws.Range("L71").Value = 36093.2715
$ws.Range("M71").Value = -6744
$ws.Range("N71").Value = -44205.2715
$ws.Range("H86").Value = 639.75
$ws.Range("I86").Value = 638.2857
$ws.Range("K86").Value = 1914.8571
$ws.Range("M86").Value = -728.8571000000002
$ws.Range("H89").Value = 639.75
$ws.Range("I89").Value = 638.2857
$ws.Range("K89").Value = 5744.571300000001
$ws.Range("M89").Value = 183.4286999999995
$ws.Range("H116").Value = 16784044
$ws.Range("I116").Value = 25174442
$ws.Range("K116").Value = 75523326
$ws.Range("M116").Value = -75519884
$ws.Range("H122").Value = 13454340
$ws.Range("J122").Value = 4730559
$ws.Range("L122").Value = 42575031
$ws.Range("N122").Value = -42579931
$ws.Range("H137").Value = 5454.6
$ws.Range("I137").Value = 2849.9092
$ws.Range("J137").Value = 7501.143
$ws.Range("K137").Value = 8549.7276
$ws.Range("L137").Value = 22503.429
$ws.Range("M137").Value = -3449.7276
$ws.Range("N137").Value = -32703.429

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 2950.2
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 2950.2
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 2950.2
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = -4008.2
$ws.Range("H62").Value = 21999.75
$ws.Range("J62").Value = 22999.666
$ws.Range("L62").Value = 22999.666
$ws.Range("N62").Value = -24371.666
$ws.Range("H65").Value = 21999.75
$ws.Range("J65").Value = 22999.666
$ws.Range("L65").Value = 68998.99800000001
$ws.Range("N65").Value = -75862.99800000001
$ws.Range("H107").Value = 801.7
$ws.Range("I107").Value = 612.8
$ws.Range("J107").Value = 990.6
$ws.Range("K107").Value = 612.8
$ws.Range("L107").Value = 990.6
$ws.Range("M107").Value = 1307.2
$ws.Range("N107").Value = -4830.6
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H122").Value = 5403.154
$ws.Range("I122").Value = 3188.95
$ws.Range("K122").Value = 9566.849999999999
$ws.Range("M122").Value = -7116.849999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5296.3335
$ws.Range("I61").Value = 4495.8184
$ws.Range("J61").Value = 7497.75
$ws.Range("K61").Value = 4495.8184
$ws.Range("L61").Value = 7497.75
$ws.Range("M61").Value = -4293.8184
$ws.Range("N61").Value = -7901.75
$ws.Range("H113").Value = 5296.3335
$ws.Range("I113").Value = 4495.8184
$ws.Range("J113").Value = 7497.75
$ws.Range("K113").Value = 4495.8184
$ws.Range("L113").Value = 7497.75
$ws.Range("M113").Value = -2325.8184
$ws.Range("N113").Value = -11837.75
$ws.Range("H122").Value = 6475.057
$ws.Range("J122").Value = 9446.799999999999
$ws.Range("L122").Value = 28340.4
$ws.Range("N122").Value = -33240.39999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2055.4443
$ws.Range("I14").Value = 1250
$ws.Range("J14").Value = 3666.3333
$ws.Range("K14").Value = 1250
$ws.Range("L14").Value = 3666.3333
$ws.Range("M14").Value = -1082
$ws.Range("N14").Value = -4002.3333
$ws.Range("H81").Value = 6333
$ws.Range("I81").Value = 3999.5
$ws.Range("K81").Value = 7999
$ws.Range("M81").Value = -6938
$ws.Range("H84").Value = 6333
$ws.Range("I84").Value = 3999.5
$ws.Range("K84").Value = 39995
$ws.Range("M84").Value = -34691
$ws.Range("H100").Value = 1432.8572
$ws.Range("I100").Value = 632
$ws.Range("J100").Value = 2500.6667
$ws.Range("K100").Value = 1264
$ws.Range("L100").Value = 5001.3334
$ws.Range("M100").Value = -723
$ws.Range("N100").Value = -6083.3334
$ws.Range("H107").Value = 6226.4443
$ws.Range("I107").Value = 607
$ws.Range("K107").Value = 1821
$ws.Range("M107").Value = 99
$ws.Range("H113").Value = 5816
$ws.Range("I113").Value = 11136.8
$ws.Range("J113").Value = 1382
$ws.Range("K113").Value = 33410.39999999999
$ws.Range("L113").Value = 4146
$ws.Range("M113").Value = -31240.39999999999
$ws.Range("N113").Value = -8486
$ws.Range("H132").Value = 1633059.6
$ws.Range("J132").Value = 7334014
$ws.Range("L132").Value = 22002042
$ws.Range("N132").Value = -22007102
$ws.Range("H136").Value = 527140.9399999999
$ws.Range("I136").Value = 1388
$ws.Range("K136").Value = 4164
$ws.Range("M136").Value = -1614
